$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add Sheet2 right after Sheet1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 with brand names, bold like Sheet1
$ws2.Range("A1").Value = "Dettol"
$ws2.Range("A2").Value = "Savlon"
$ws2.Range("A1:A2").Font.Bold = $true

# Match the selection/active cell recorded in the target sheet
$null = $ws2.Range("D3").Select()
